# Weekly update: insert a new "Arveja Verde" price record for
# Terminal Hortofrutícola Agro Chillán, shifting the existing rows
# 59-66 down to 60-67 and adding a new row 59 with this week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 59 (old rows 59-66 -> 60-67)
$ws.Rows.Item(59).Insert()

# Populate the newly inserted row 59 with this week's record
$ws.Cells.Item(59, 1).Value  = 7
$ws.Cells.Item(59, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(59, 3).Value  = "Ñuble"
$ws.Cells.Item(59, 4).Value  = 44615
$ws.Cells.Item(59, 5).Value  = 16
$ws.Cells.Item(59, 6).Value  = 100112022
$ws.Cells.Item(59, 7).Value  = "Arveja Verde"
$ws.Cells.Item(59, 8).Value  = "Sin especificar"
$ws.Cells.Item(59, 9).Value  = "Primera"
$ws.Cells.Item(59, 10).Value = 60
$ws.Cells.Item(59, 11).Value = 24000
$ws.Cells.Item(59, 12).Value = 25000
$ws.Cells.Item(59, 13).Value = 24500
$ws.Cells.Item(59, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(59, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(59, 16).Value = 980
$ws.Cells.Item(59, 17).Value = 25
$ws.Cells.Item(59, 18).Value = "Hortaliza"
